# Add 60% AMI affordability information to the income affordability sheet.
#
# Before:  Year | MedianIncome | 80AMI | Affordability | MedianSale | 60AMI | RentAffordability | 3BRAveRent
# After:   Year | MedianIncome | 80AMI | 80Affordability | MedianSale | 60AMI | 60Affordability | RentAffordability | 3BRAveRent
#
# i.e. the old generic "Affordability" header (column D) is renamed to
# "80Affordability" (clarifying that it goes with the 80AMI column), and a
# brand-new "60Affordability" column is inserted right after 60AMI, holding
# the new data point (190400) that goes with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "Affordability" header (column D) to "80Affordability".
$ws.Range("D1").Value = "80Affordability"

# Insert a new column before the old "RentAffordability" column (G) for the
# new "60Affordability" data, pushing RentAffordability/3BRAveRent one
# column to the right.
$ws.Range("G1").EntireColumn.Insert()
$ws.Range("G1").Value = "60Affordability"
$ws.Range("G2").Value = 190400

# Match the author's final selection noted in the saved file.
$ws.Range("G10").Select()
